# JT -add unity game
#
# The underlying text is unchanged in both spots touched by this edit; what
# moved is Word's automatic "last edit" bookmark (_GoBack). Word keeps only
# one _GoBack bookmark in the whole document and relocates it to wherever
# the most recent edit happened, splitting the run(s) at that point as
# needed.
#
# 1) The previous edit location - inside "Documented technical in|structions"
#    - is no longer the latest edit, so that bookmark goes away and the run
#    it used to split is rejoined into a single run.
# 2) The new edit happened right after "...competitive position", so
#    _GoBack now sits there (as an empty/collapsed bookmark) and the run
#    that used to span "ecure a challenging ... professional growth." is
#    split in two around that point.

$d = $word.ActiveDocument

# --- Rejoin "Documented technical in" + "structions" into a single run,
#     which also removes the stale _GoBack bookmark that used to separate
#     them (Find/Replace across the bookmarked split re-merges the runs).
$d.Content.Find.Execute("Documented technical instructions", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Documented technical instructions", 2) | Out-Null

# --- Move _GoBack to the new edit point: right after "...competitive
#     position", before " with the opportunity for professional growth."
$rng = $d.Content
$rng.Find.Execute("ecure a challenging and competitive position", $true, `
    $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$editPoint = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $editPoint) | Out-Null
